$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per Dr Hou advice: ligand/receptor-expressing cell counts (columns E and K)
# increase from 1 to 3 for every data row, and all derived expression/
# specificity metrics (columns G-J, M-T) are recomputed accordingly.

$rowData = @{}
$rowData[2] = @(3, 1, 0.135131, 0.405393, 0.1435196772981866, 0.1435196772981866, 3, 1, 144.283282, 432.849846, 0.5695387292414458, 0.5695387292414458, 19.497144179942, 175.474297619478, 0.08174001462955155, 0.08174001462955155)
$rowData[3] = @(3, 1, 0.135131, 0.405393, 0.1435196772981866, 0.1435196772981866, 3, 1, 29.99441333333333, 89.98324, 0.118398887365269, 0.118398887365269, 4.053175068146667, 36.47857561332, 0.01699257010712775, 0.01699257010712775)
$rowData[4] = @(3, 1, 0.135131, 0.405393, 0.1435196772981866, 0.1435196772981866, 3, 1, 31.60894133333333, 94.826824, 0.1247720181445143, 0.1247720181445143, 4.271347851314667, 38.442130661832, 0.01790723977994417, 0.01790723977994417)
$rowData[5] = @(3, 1, 0.135131, 0.405393, 0.1435196772981866, 0.1435196772981866, 3, 1, 32.77231833333333, 98.31695499999999, 0.1293642914073912, 0.1293642914073912, 4.428556148701666, 39.857005338315, 0.01856632135669736, 0.01856632135669736)
$rowData[6] = @(3, 1, 0.135131, 0.405393, 0.1435196772981866, 0.1435196772981866, 3, 1, 14.67461933333333, 44.023858, 0.05792607384137975, 0.05792607384137975, 1.982995985132667, 17.846963866194, 0.008313531424865749, 0.008313531424865749)
$rowData[7] = @(3, 1, 0.8064193333333334, 2.419258, 0.8564803227018134, 0.8564803227018134, 3, 1, 144.283282, 432.849846, 0.5695387292414458, 0.5695387292414458, 116.3528280815854, 1047.175452734268, 0.4877987146118942, 0.4877987146118942)
$rowData[8] = @(3, 1, 0.8064193333333334, 2.419258, 0.8564803227018134, 0.8564803227018134, 3, 1, 29.99441333333333, 89.98324, 0.118398887365269, 0.118398887365269, 24.18807480399111, 217.69267323592, 0.1014063172581412, 0.1014063172581412)
$rowData[9] = @(3, 1, 0.8064193333333334, 2.419258, 0.8564803227018134, 0.8564803227018134, 3, 1, 31.60894133333333, 94.826824, 0.1247720181445143, 0.1247720181445143, 25.49006139739911, 229.410552576592, 0.1068647783645701, 0.1068647783645701)
$rowData[10] = @(3, 1, 0.8064193333333334, 2.419258, 0.8564803227018134, 0.8564803227018134, 3, 1, 32.77231833333333, 98.31695499999999, 0.1293642914073912, 0.1293642914073912, 26.42823110215445, 237.85407991939, 0.1107979700506938, 0.1107979700506938)
$rowData[11] = @(3, 1, 0.8064193333333334, 2.419258, 0.8564803227018134, 0.8564803227018134, 3, 1, 14.67461933333333, 44.023858, 0.05792607384137975, 0.05792607384137975, 11.83389673970711, 106.505070657364, 0.049612542416514, 0.049612542416514)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
